$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the old C21 formula up to C19 (row 21 -> row 19) ---
$ws.Range("C21").ClearContents()
$ws.Range("C19").Formula = "=(5/(2*PI()*D22))"

# --- New label cells around the LINEST results (rows 22-23) ---
# Order matters for shared-string table insertion order:
# gradient(19), ∆gradient(20), intercept(21), ∆intercept(22)
$ws.Range("C22").Value = "gradient"
$ws.Range("C23").Value = "∆gradient"
$ws.Range("F22").Value = "intercept"
$ws.Range("F23").Value = "∆intercept"

# Apply the "body" font style (non-theme font) used for the ∆ row labels
$ws.Range("C23").Font.Name = "Calibri"
$ws.Range("F23").Font.Name = "Calibri"

# --- New summary table restating gradient/intercept with uncertainties ---
$ws.Range("H25").Value = "gradient"
$ws.Range("I25").Value = 2.432703935776543
$ws.Range("J25").Formula = "=-7.52532*10^-5"
$ws.Range("K25").Value = "intercept"

$ws.Range("H26").Value = "∆gradient"
$ws.Range("I26").Value = 0.04482739828602023
$ws.Range("J26").Value = 0.000018241036805583528
$ws.Range("K26").Value = "∆intercept"

# --- Move the percentage-error formula from D25 down to D27 ---
$ws.Range("D25").ClearContents()
$ws.Range("D27").Formula = "=(D23/D22)*100"

# --- Column J width tweak (slightly narrower) ---
$ws.Columns("J").ColumnWidth = 10.5

# --- View state: scrolled down a couple more rows, new active selection ---
$ws.Range("L28").Select() | Out-Null
